# "playing with physics calcs" - tweak the projectile Time input and
# re-point the UI selection, matching the author's interactive session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Time (F2) goes from 0.85s to 1.5s; every dependent formula in
# H:I / K:L / N:O recalculates automatically from this single input.
$ws.Range("F2").Value = 1.5

# Leave the cursor on F3, where the author was last working.
$ws.Range("F3").Select()

# Reposition the saved window chrome to match the author's session.
$win = $wb.Windows.Item(1)
$win.Left = 2400
$win.Top = 2100
